# ---------------------------------------------------------------------------
# read dtc info report by status mask implemented
#
# This script reproduces (as closely as the COM surface allows) the edits
# described by the target diff:
#   - New ReadDtcInfoGeneric / ReadDtcInfoGenericPos / ReadDTCInfoBySt /
#     ReadDTCInfoByStPos / DtcNStatusRecordPos rows appended to the
#     "ReqResp" worksheet (rows 60-76), which also introduces 15 new
#     shared strings and grows the XML-mapped table (Table8) accordingly.
#   - The helper "Req"/"Resp" list used by the data-validation drop-down
#     (originally parked in column N) is moved one column to the left
#     (column M) after the now-unused column M is removed.
#   - Minor selection / cursor position bookkeeping on a few sheets.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Services sheet: just a cursor position change left over from editing.
# ---------------------------------------------------------------------------
$wsServices = $wb.Worksheets.Item("Services")
$wsServices.Activate()
$wsServices.Range("B20").Select()

# ---------------------------------------------------------------------------
# 2. ServiceSub sheet: cursor position / scroll bookkeeping.
# ---------------------------------------------------------------------------
$wsServiceSub = $wb.Worksheets.Item("ServiceSub")
$wsServiceSub.Activate()
$wsServiceSub.Range("B116").Select()

# ---------------------------------------------------------------------------
# 3. ReqResp sheet: the actual data/content change.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ReqResp")
$ws.Activate()

# --- 3a. Append the new ReqResp rows for ReadDtcInfoByStatusMask support ---

# ReadDtcInfoGeneric (request)
$ws.Range("A60").Value = "ReadDtcInfoGeneric"
$ws.Range("B60").Value = "sid"
$ws.Range("C60").Value = "Req"
$ws.Range("D60").Value = 1
$ws.Range("E60").Value = 1

$ws.Range("A61").Value = "ReadDtcInfoGeneric"
$ws.Range("B61").Value = "subFunc"
$ws.Range("C61").Value = "Req"
$ws.Range("D61").Value = 2
$ws.Range("E61").Value = 1

$ws.Range("A62").Value = "ReadDtcInfoGeneric"
$ws.Range("B62").Value = "buf"
$ws.Range("C62").Value = "Req"
$ws.Range("D62").Value = 3
$ws.Range("E62").Value = "n"

# ReadDtcInfoGenericPos (positive response)
$ws.Range("A63").Value = "ReadDtcInfoGenericPos"
$ws.Range("B63").Value = "sid"
$ws.Range("C63").Value = "Resp"
$ws.Range("D63").Value = 1
$ws.Range("E63").Value = 1

$ws.Range("A64").Value = "ReadDtcInfoGenericPos"
$ws.Range("B64").Value = "subFunc"
$ws.Range("C64").Value = "Resp"
$ws.Range("D64").Value = 2
$ws.Range("E64").Value = 1

$ws.Range("A65").Value = "ReadDtcInfoGenericPos"
$ws.Range("B65").Value = "buf"
$ws.Range("C65").Value = "Resp"
$ws.Range("D65").Value = 3
$ws.Range("E65").Value = "n"

# ReadDTCInfoBySt (request)
$ws.Range("A66").Value = "ReadDTCInfoBySt"
$ws.Range("B66").Value = "sid"
$ws.Range("C66").Value = "Req"
$ws.Range("D66").Value = 1
$ws.Range("E66").Value = 1

$ws.Range("A67").Value = "ReadDTCInfoBySt"
$ws.Range("B67").Value = "subFunc"
$ws.Range("C67").Value = "Req"
$ws.Range("D67").Value = 2
$ws.Range("E67").Value = 1

$ws.Range("A68").Value = "ReadDTCInfoBySt"
$ws.Range("B68").Value = "mask"
$ws.Range("C68").Value = "Req"
$ws.Range("D68").Value = 3
$ws.Range("E68").Value = 1

# ReadDTCInfoByStPos (positive response)
$ws.Range("A69").Value = "ReadDTCInfoByStPos"
$ws.Range("B69").Value = "sid"
$ws.Range("C69").Value = "Resp"
$ws.Range("D69").Value = 1
$ws.Range("E69").Value = 1

$ws.Range("A70").Value = "ReadDTCInfoByStPos"
$ws.Range("B70").Value = "reportType"
$ws.Range("C70").Value = "Resp"
$ws.Range("D70").Value = 2
$ws.Range("E70").Value = 1

$ws.Range("A71").Value = "ReadDTCInfoByStPos"
$ws.Range("B71").Value = "availStMask"
$ws.Range("C71").Value = "Resp"
$ws.Range("D71").Value = 3
$ws.Range("E71").Value = 1

$ws.Range("A72").Value = "ReadDTCInfoByStPos"
$ws.Range("B72").Value = "dtcNStRecord"
$ws.Range("C72").Value = "Resp"
$ws.Range("D72").Value = 4
$ws.Range("E72").Value = "n"

# DtcNStatusRecordPos (content of dtcNStRecord)
$ws.Range("A73").Value = "DtcNStatusRecordPos"
$ws.Range("B73").Value = "dtcHighByte"
$ws.Range("C73").Value = "Resp"
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 1
$ws.Range("F73").Value = "ReadDTCInfoByStPos dtcNStRecord content"

$ws.Range("A74").Value = "DtcNStatusRecordPos"
$ws.Range("B74").Value = "dtcMiddleByte"
$ws.Range("C74").Value = "Resp"
$ws.Range("D74").Value = 2
$ws.Range("E74").Value = 1

$ws.Range("A75").Value = "DtcNStatusRecordPos"
$ws.Range("B75").Value = "dtcLowByte"
$ws.Range("C75").Value = "Resp"
$ws.Range("D75").Value = 3
$ws.Range("E75").Value = 1

$ws.Range("A76").Value = "DtcNStatusRecordPos"
$ws.Range("B76").Value = "statusOfDtc"
$ws.Range("C76").Value = "Resp"
$ws.Range("D76").Value = 4
$ws.Range("E76").Value = 1

# --- 3b. Re-apply the "group name" wrap/centre style to column A for the ---
# --- two new multi-row groups (matches the style already used elsewhere) --
$groupA = $ws.Range("A66:A72")
$groupA.NumberFormat = "@"
$groupA.VerticalAlignment = -4108
$groupA.WrapText = $true

$ws.Range("F73").NumberFormat = "@"
$ws.Range("F73").VerticalAlignment = -4108
$ws.Range("F73").WrapText = $true

# --- 3c. Grow the XML-mapped table (Table8) so it covers the new rows. ----
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F76"))

# --- 3d. The helper Req/Resp list used for the column-C drop-down lived --
# --- in column N (empty column M in between); remove the empty column ---
# --- M so the helper values end up in M instead of N, and fix up the ----
# --- data validation formula to match. -----------------------------------
$ws.Columns.Item(13).Delete()

$ws.Range("C1:C1048576").Validation.Delete()
$ws.Range("C1:C1048576").Validation.Add(3, 1, 1, "=`$M`$5:`$M`$6")
$ws.Range("C1:C1048576").Validation.IgnoreBlank = $true
$ws.Range("C1:C1048576").Validation.InCellDropdown = $true
$ws.Range("C1:C1048576").Validation.ShowInput = $true
$ws.Range("C1:C1048576").Validation.ShowError = $true
$ws.Range("C8").Validation.Delete()

# --- 3e. Column A got a bit wider once the longer names were added. -----
$ws.Columns.Item(1).ColumnWidth = 24.7142857142857

# --- 3f. Cursor / selection bookkeeping. ---------------------------------
$ws.Range("C63:C65").Select()
